$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 13.377
$ws.Range("B10").Value = 7.095000000000001
$ws.Range("B12").Value = 6.548999999999999
$ws.Range("D13").Value = -7.831999999999999
$ws.Range("B18").Value = 6.548999999999999
$ws.Range("E20").Value = 12.932
